$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 2) text changes
$ws.Range("C2").Value = "dense embedding  generálásai idő átlaga"
$ws.Range("G2").Value = "teljes feldoldozási idő átlaga (back-end)"

# Row 3 (20 questions) values
$ws.Range("D3").Value = "2.47/s"
$ws.Range("E3").Value = "0.62/s"
$ws.Range("F3").Value = "2.46/s"
$ws.Range("G3").Value = "6.03/s"
$ws.Range("H3").Value = 0.89

# Row 4 (40 questions) values
$ws.Range("C4").Value = "0.46/s"
$ws.Range("D4").Value = "2.44/s"
$ws.Range("E4").Value = "0.61/s"
$ws.Range("F4").Value = "1.65/s"
$ws.Range("G4").Value = "5.17/s"
$ws.Range("H4").Value = 0.84

# Row 5 (60 questions) values
$ws.Range("C5").Value = "0.46/s"
$ws.Range("D5").Value = "2.44/s"
$ws.Range("F5").Value = "1.66/s"
$ws.Range("G5").Value = "5.17/s"

# Update the selection to match the new saved view state
$ws.Range("F7").Select()
